# Big changes. Added the repeater functions for a potential all night
# experiment. Ran trades.
#
# The repeater ran three new trades (rows 3-5) with the same columns as
# the existing row 2, and the original row 2 was backfilled with the
# newly-computed PriceChange (X) / UpDown (Y) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
$ws.Range("A3").Value = 42649.671851851854
$ws.Range("A3").NumberFormat = "m/d/yy h:mm"
$ws.Range("B3").Value = -17
$ws.Range("C3").Value = "Strong Sell"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = "Random"
$ws.Range("Q3").Value = 35.483819709395277
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0.0965
$ws.Range("S3").NumberFormat = "0.00%"
$ws.Range("T3").Value = 0.0269
$ws.Range("T3").NumberFormat = "0.00%"
$ws.Range("U3").Value = 4.82
$ws.Range("V3").Value = 2.28
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = -0.35000000000000853
$ws.Range("Y3").Value = "Down"

# --- Row 4 ---
$ws.Range("A4").Value = 42649.672175925924
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "Neutral"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = "Random"
$ws.Range("Q4").Value = 35.483819709395277
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0.0965
$ws.Range("S4").NumberFormat = "0.00%"
$ws.Range("T4").Value = 0.0269
$ws.Range("T4").NumberFormat = "0.00%"
$ws.Range("U4").Value = 4.82
$ws.Range("V4").Value = 2.28
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = -0.35000000000000853
$ws.Range("Y4").Value = "Down"

# --- Row 5 ---
$ws.Range("A5").Value = 42649.674039351848
$ws.Range("A5").NumberFormat = "m/d/yy h:mm"
$ws.Range("B5").Value = -5
$ws.Range("C5").Value = "Sell"
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = "Random"
$ws.Range("Q5").Value = 35.483819709395277
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0.0965
$ws.Range("S5").NumberFormat = "0.00%"
$ws.Range("T5").Value = 0.0269
$ws.Range("T5").NumberFormat = "0.00%"
$ws.Range("U5").Value = 4.82
$ws.Range("V5").Value = 2.28
$ws.Range("W5").Value = 0

# --- Backfill row 2 with the newly-computed PriceChange / UpDown columns ---
$ws.Range("X2").Value = -0.35000000000000853
$ws.Range("Y2").Value = "Down"

# Column C ("Verdict") needs to widen now that "Strong Sell" lives in it.
$ws.Columns.Item(3).ColumnWidth = 8.75
